$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-05 Friday" "2025-09-06 Saturday"

Replace-Text "13÷9=1, 4" "40÷6=6, 4"
Replace-Text "44÷7=6, 2" "78÷9=8, 6"
Replace-Text "85÷5=17, 0" "51÷8=6, 3"
Replace-Text "42÷6=7, 0" "38÷4=9, 2"
Replace-Text "75÷2=37, 1" "50÷4=12, 2"

Replace-Text "19÷6=3, 1" "37÷5=7, 2"
Replace-Text "17÷5=3, 2" "94÷7=13, 3"
Replace-Text "76÷6=12, 4" "44÷4=11, 0"
Replace-Text "95÷8=11, 7" "54÷8=6, 6"
Replace-Text "81÷2=40, 1" "96÷2=48, 0"

Replace-Text "54÷7=7, 5" "88÷9=9, 7"
Replace-Text "70÷9=7, 7" "84÷2=42, 0"
Replace-Text "68÷4=17, 0" "63÷4=15, 3"
Replace-Text "42÷2=21, 0" "32÷5=6, 2"
Replace-Text "23÷5=4, 3" "61÷7=8, 5"

Replace-Text "74÷9=8, 2" "27÷5=5, 2"
Replace-Text "84÷7=12, 0" "68÷7=9, 5"
Replace-Text "14÷3=4, 2" "65÷6=10, 5"
Replace-Text "22÷6=3, 4" "34÷9=3, 7"
Replace-Text "64÷8=8, 0" "20÷6=3, 2"

Replace-Text "78÷2=39, 0" "84÷3=28, 0"
Replace-Text "35÷9=3, 8" "54÷4=13, 2"
Replace-Text "80÷4=20, 0" "94÷7=13, 3"
Replace-Text "42÷9=4, 6" "79÷8=9, 7"
Replace-Text "14÷7=2, 0" "20÷8=2, 4"

Write-Output "done"
